$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: Date (as text, matching existing shared-string date cells) and USDValue
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2024-11-15"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = 0.03313
